# Update invoice report column order and widths
# Reorders columns C..K on the "GST Report" sheet and adjusts a few column
# widths so the layout matches the newly requested report format:
#   old: S.No | Vendor | GSTIN | Invoice No. | Date | Taxable Amt | Total Tax | CGST | SGST | IGST | HSN Codes
#   new: S.No | Vendor | Date  | GSTIN | Invoice No. | HSN Codes | CGST | SGST | IGST | Total Tax | Taxable Amt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 2   # header (row 1) + a single data row (row 2) in this report
$lastCol = 11  # columns A..K

# ---------------------------------------------------------------------
# 1) Re-write the header row (row 1) in the new column order.
#    Header cells are plain text, so a direct .Value assignment is safe.
# ---------------------------------------------------------------------
$newHeaders = @{
    1  = "S.No."
    2  = "Vendor/Shop Name"
    3  = "Date"
    4  = "GSTIN"
    5  = "Invoice No."
    6  = "HSN Codes"
    7  = "CGST"
    8  = "SGST"
    9  = "IGST"
    10 = "Total Tax"
    11 = "Taxable Amount"
}
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item(1, $col).Value = $newHeaders[$col]
}

# ---------------------------------------------------------------------
# 2) Re-order the data row (row 2) to match the new header order.
#
#    Several of the values look like numbers/dates (e.g. "2024-09-28",
#    "173.91", "4045.01") but must stay plain text, exactly as they were
#    authored. Re-typing them through .Value would let Excel "helpfully"
#    re-interpret them as a real date/number. To avoid that (and to avoid
#    touching NumberFormat, which creates unwanted extra cell styles) we
#    stage the original row in a scratch area, then copy the values back
#    into their new homes with PasteSpecial (values only). Copying a
#    cell that is already text-typed keeps it text-typed, and
#    PasteSpecial-values never touches the destination cell's existing
#    style/format.
# ---------------------------------------------------------------------
$stageRow = 500
$srcRange   = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
$stageRange = $ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, $lastCol))
$srcRange.Copy()
$stageRange.PasteSpecial(-4163)   # xlPasteValues

# target column -> original column that its value comes from
$colSource = @{
    1  = 1    # S.No.            <- S.No.
    2  = 2    # Vendor/Shop Name <- Vendor/Shop Name
    3  = 5    # Date             <- Date
    4  = 3    # GSTIN            <- GSTIN
    5  = 4    # Invoice No.      <- Invoice No.
    6  = 11   # HSN Codes        <- HSN Codes
    7  = 8    # CGST             <- CGST
    8  = 9    # SGST             <- SGST
    9  = 10   # IGST             <- IGST
    10 = 7    # Total Tax        <- Total Tax
    11 = 6    # Taxable Amount   <- Taxable Amount
}

# Copy highest source columns first isn't required since we already staged
# a full snapshot of row 2, so destination writes can happen in any order
# without clobbering a value we still need to read.
for ($col = 1; $col -le $lastCol; $col++) {
    $srcCol = $colSource[$col]
    $from = $ws.Cells.Item($stageRow, $srcCol)
    $to   = $ws.Cells.Item(2, $col)
    $from.Copy()
    $to.PasteSpecial(-4163)       # xlPasteValues - keeps $to's own style/format
}

# Clean up the scratch area used for staging.
$stageRange.ClearContents()

# ---------------------------------------------------------------------
# 3) Swap the wrap-text formatting between the (now relocated) HSN Codes
#    column (F) and Taxable Amount column (K): the long comma separated
#    HSN code list needs to wrap, the plain numeric amount does not.
# ---------------------------------------------------------------------
$ws.Range("F2").WrapText = $true
$ws.Range("K2").WrapText = $false

# ---------------------------------------------------------------------
# 4) Apply the updated column widths.
#    NOTE: Excel's ColumnWidth property is expressed in "characters" but
#    the value actually persisted in the OOXML <col width> is offset by
#    a constant 5/6 (~0.8333) padding factor. Subtract that offset here
#    so the saved width exactly matches the requested value.
# ---------------------------------------------------------------------
$padding = 5 / 6
$newWidths = @{
    3  = 15
    4  = 18
    5  = 20
    6  = 40
    11 = 15
}
foreach ($col in $newWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $newWidths[$col] - $padding
}
